$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '30.189.58'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +0.60%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.079.02'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -1.87%  '
$ws.Range('E4').Value = '  -0.48%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '338.44'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -2.76%  '
$ws.Range('E6').Value = '  -0.54%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.5248'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.86%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.4356'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -2.27%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '54.92'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +1.20%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.09337'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -0.27%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '1.168'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -1.00%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '24.46'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -3.10%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '8.441'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.17%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.841'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.64%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '2.052.23'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -1.82%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '100.20'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -2.34%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.00001155'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -0.77%  '
$ws.Range('E18').Value = '  -0.47%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '20.81'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -3.81%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.06699'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.06%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.280'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.44%  '
$ws.Range('E22').Value = '  -0.56%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '30.221.05'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.68%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '12.36'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -3.31%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.311'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.66%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '21.71'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -1.91%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '162.20'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.36%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '6.777'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +3.28%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.477'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -2.85%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '133.11'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.82%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.125'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -2.57%  '
$ws.Range('B32').Value = 'Stellar'
$ws.Range('C32').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.1045'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -1.15%  '
$ws.Range('B33').Value = 'ARBITRUM'
$ws.Range('C33').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.655'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -7.07%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '6.234'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.35%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '3.909'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -1.68%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.02597'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.36%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '9.826'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -9.03%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.06704'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -2.48%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.6926'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.64%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '12.48'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -1.94%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.323'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.90%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.2195'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -2.45%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.6702'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -2.34%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.350'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.15%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '14.17'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -3.12%  '
$ws.Range('E46').Value = '  -0.53%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.306'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +4.89%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '3.614'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.76%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.00000000351'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -2.51%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.206'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +1.69%  '
$ws.Range('E51').Value = '  -1.63%  '
